# Apply the edits described by the commit "add Output of values to the interface"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: turbine marking / manufacturer and input parameters ---
$ws.Range("B3").Value = "К-300-23,5"
$ws.Range("C3").Value = "ЛМЗ, ЦВД"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 3
$ws.Range("I3").Value = 2
$ws.Range("K3").Value = 0

# --- Row 9: initial steam parameters (Индекс 0) ---
$ws.Range("D9").Value = 13.97868135173841
$ws.Range("E9").Value = 0.04362872542254398
$ws.Range("F9").Value = 0.0009990908150413771

# --- Row 14: final steam parameters (Индекс kt') ---
$ws.Range("C14").Value = 2.999999999999773
$ws.Range("D14").Value = 13.97868135173728
$ws.Range("F14").Value = 0.0009990908150413771

# --- Row 18: available heat drop H0 (value uses scientific notation in source) ---
$ws.Range("C18").Value = 0.00000000000113509202037676

# --- Row 22: nominal consumption G0 (value uses scientific notation in source) ---
$ws.Range("C22").Value = 12437447189645120

# --- Row 24: electrical power Nэ ---
$ws.Range("H24").Value = 2027.239756462033

# --- Row 27: pressure drop estimate p0' ---
$ws.Range("F27").Value = 1.9

# --- Row 32: steam parameters after control valves (Индекс 0 (штрих)) ---
$ws.Range("C32").Value = 3.023854003193946
$ws.Range("D32").Value = 13.97868135173841
$ws.Range("E32").Value = 0.04399070813284934
$ws.Range("F32").Value = 0.0009991400950681406
